$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values look like plain numbers,
# so Excel keeps them as text (matching the original inlineStr cell type) instead of
# silently coercing them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated price / volume values.
$ws.Range("D2").Value = "67.872.61"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "2.403.76"
$ws.Range("E3").Value = "  -3.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "553.54"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "158.27"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "0.162"
$ws.Range("E9").Value = "  +5.34%  "
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "4.64"
$ws.Range("E12").Value = "  -5.07%  "
$ws.Range("D13").Value = "67.759.34"
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("D14").Value = "2.845.38"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "22.74"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").Value = "2.407.62"
$ws.Range("E17").Value = "  -3.02%  "
$ws.Range("D18").Value = "10.31"
$ws.Range("E18").Value = "  -5.08%  "
$ws.Range("D19").Value = "329.44"
$ws.Range("E19").Value = "  -3.36%  "
$ws.Range("D20").Value = "6.83"
$ws.Range("E20").Value = "  -3.79%  "
$ws.Range("D21").Value = "3.77"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  -3.90%  "
$ws.Range("D24").Value = "65.92"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").Value = "3.63"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").Value = "2.532.09"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "8.12"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "0.0₃0804"
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("D29").Value = "7.05"
$ws.Range("E29").Value = "  -2.78%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "416.71"
$ws.Range("E31").Value = "  -4.83%  "
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").Value = "1.59"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").Value = "159.23"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "18.98"
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("D37").Value = "17.72"
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("D39").Value = "0.294"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").Value = "4.27"
$ws.Range("E40").Value = "  -4.43%  "
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").Value = "1.06"
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("D43").Value = "131.12"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D44").Value = "3.29"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("D45").Value = "1.95"
$ws.Range("E45").Value = "  -7.11%  "
$ws.Range("D46").Value = "0.0708"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").Value = "0.476"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "0.552"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").Value = "0.0911"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "1.36"
$ws.Range("E51").Value = "  -4.67%  "
